$wb = $excel.ActiveWorkbook

$wsMeta = $wb.Worksheets.Item("openbis-metadata")
$wsData = $wb.Worksheets.Item("openbis-data")

# Replace the (previously valid) CHEBI compound identifiers with bad/placeholder
# data so the metabolomics *data* rows get exercised by validation too, not
# just the metadata sheet.
$wsData.Range("A2").Value = "foo"
$wsData.Range("A3").Value = "foo"

# Move the active selection / active tab: the metadata sheet keeps a stale
# selection, and the data sheet becomes the active tab with its own
# selection.
$wsMeta.Range("D15").Select()
$wsData.Activate()
$wsData.Range("A4").Select()
